$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data as produced by the scheduled scraper run.
# Setting NumberFormat to text ("@") before assigning values preserves the original
# text representation (e.g. "28.911.39", "0.9995") instead of letting Excel coerce
# numeric-looking strings into floating point numbers.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.911.39"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.76%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.832.50"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.91%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.22%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.28"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.63%  "

# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.31%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9997"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.22%  "

# Row 8
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07652"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.89%  "

# Row 9
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3053"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.56%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.50"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.23%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07817"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.40%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.833.22"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.55%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.064"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.40%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "90.35"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6783"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.22%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.418"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.35%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008316"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.46%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "28.916.03"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.03%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.58"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.81%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.082.34"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.27%  "

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.32%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9992"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.24%  "

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.52%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9996"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.21%  "

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -5.67%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.92"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.49%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.796"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.28%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.19"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.84%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.559"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.68%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.216"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.66%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.141"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.72%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.177"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.26%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05114"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.84%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7554"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.04%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.842"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.71%  "

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.23%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.678"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.36%  "

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.51%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.227.59"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.03%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.689"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.84%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9270"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.57%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "108.85"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.14%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9991"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.17%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.718"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -6.14%  "

# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.981.53"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.46%  "

# Row 46
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5172"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.28%  "

# Row 47
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.534"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.07%  "

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.93%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "64.12"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -10.26%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.738"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.10%  "

# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.909"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.98%  "
